$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove all existing hyperlinks up front -------------------------------
# (Range-scoped Hyperlinks.Delete() on this host always clears the whole
# worksheet collection, so we reset once and rebuild every hyperlink below
# in the order the target workbook expects.)
$ws.Range("A1").Hyperlinks.Delete()

# --- Insert the new row 13 --------------------------------------------------
# This pushes the current row 13 (DRA051 / OPQA-4188||OPQA-4190 / ...) down
# to row 14, carrying its own cell formatting with it.
$ws.Rows.Item(13).Insert()

# --- Populate the new row 13 with the OPQA-4223 / OPQA-4224 test case ------
$ws.Range("A13").Value = "DRA006"
$ws.Range("B13").Value = "OPQA-4223 || OPQA-4224"
$ws.Range("C13").Value = 'Verify that error message " New password should not match current password" should be displayed when user enters the current password in change password field.|| Verify that error message"New password should not match previous 4 passwords" should be displayed when user enters password in change password field which is matching with the previous 4 passwords.'
$ws.Range("D13").Value = "Y"

# --- Row heights -------------------------------------------------------------
$ws.Rows.Item(13).RowHeight = 75
$ws.Rows.Item(14).RowHeight = 60

# --- Rebuild every hyperlink, in the same order as the edited workbook -----
$ws.Hyperlinks.Add($ws.Range("B2"), "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-4176", "", "", "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-4176")
$ws.Range("B2").Value = "OPQA-4176||OPQA-4178||OPQA-4179||OPQA-4182||OPQA-4187||OPQA-4189 "

$ws.Hyperlinks.Add($ws.Range("B14"), "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-4221", "", "", "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-4221")
$ws.Range("B14").Value = "OPQA-4188||OPQA-4190"

$ws.Hyperlinks.Add($ws.Range("B12"), "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-4221", "", "", "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-4221")
$ws.Range("B12").Value = "OPQA-4221"

$ws.Hyperlinks.Add($ws.Range("B13"), "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-4223", "", "", "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-4223")
$ws.Range("B13").Value = "OPQA-4223 || OPQA-4224"

# --- Leave the selection where the author's last click landed --------------
$ws.Range("B13").Select()
